# Reexecution of failed testcases logic implementation
#
# Insert a new row (CLICK_PRE_ENTERTEXT) above the existing ENTERTEXT row on
# the first worksheet, shifting the remaining rows down by one, then update
# the sheet selection to reflect the new cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 3 (pushes old rows 3-6 down to 4-7)
$ws.Rows.Item(3).Insert()

# Give the new row the same bordered look as the surrounding data rows
$newRowRange = $ws.Range("A3:E3")
$newRowRange.Borders.LineStyle = 1
$newRowRange.Borders.Weight = 2

# Populate the new row: CLICK_PRE_ENTERTEXT / SearchBoxHomePage / CSS / (blank)
$ws.Cells.Item(3, 1).Value = ""
$ws.Cells.Item(3, 2).Value = "CLICK_PRE_ENTERTEXT"
$ws.Cells.Item(3, 3).Value = "SearchBoxHomePage"
$ws.Cells.Item(3, 4).Value = "CSS"
$ws.Cells.Item(3, 5).Value = ""

# Move the visible selection to C3:D3 with C3 as the active cell
$ws.Range("C3:D3").Select()
